$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 0
    5  = -2
    6  = 2
    7  = 3
    8  = 9
    9  = 1
    10 = 5
    11 = 4
    12 = 1
    15 = 6
    16 = 2
    18 = 2
    19 = -6
    21 = 1
    22 = -5
    23 = -2
    24 = 4
    25 = -1
    26 = 2
    28 = 4
    29 = 5
    30 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
